# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.130.76"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.668.33"
$ws.Range("E3").Value = "  -1.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.75"
$ws.Range("E5").Value = "  -3.26%  "

# Row 6
$ws.Range("E6").Value = "  -2.79%  "

# Row 7
$ws.Range("E7").Value = "  -0.52%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2634"
$ws.Range("E8").Value = "  -3.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06303"
$ws.Range("E9").Value = "  -2.34%  "

# Row 10
$ws.Range("E10").Value = "  -2.50%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07547"
$ws.Range("E11").Value = "  -1.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.670.81"
$ws.Range("E12").Value = "  -1.75%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.443"
$ws.Range("E13").Value = "  -2.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5561"
$ws.Range("E14").Value = "  -3.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.81"
$ws.Range("E15").Value = "  -0.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007938"
$ws.Range("E16").Value = "  -5.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.157.50"

# Row 18
$ws.Range("E18").Value = "  -0.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.749"
$ws.Range("E19").Value = "  -3.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.58"
$ws.Range("E20").Value = "  -2.08%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.36"
$ws.Range("E21").Value = "  -4.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.178"

# Row 23
$ws.Range("E23").Value = "  -0.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.69"
$ws.Range("E24").Value = "  +0.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1251"
$ws.Range("E25").Value = "  -2.78%  "

# Row 26
$ws.Range("E26").Value = "  -4.60%  "

# Row 27
$ws.Range("E27").Value = "  +0.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06245"
$ws.Range("E28").Value = "  -0.62%  "

# Row 29
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").Value = "  -3.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.517"
$ws.Range("E31").Value = "  -2.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.419"
$ws.Range("E32").Value = "  -4.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.630"
$ws.Range("E33").Value = "  -2.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9981"
$ws.Range("E34").Value = "  -3.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6053"
$ws.Range("E35").Value = "  -1.95%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.413"
$ws.Range("E36").Value = "  -0.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  -1.49%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.136"
$ws.Range("E38").Value = "  +0.28%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.104.99"
$ws.Range("E39").Value = "  -0.55%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01615"
$ws.Range("E40").Value = "  -2.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8720"
$ws.Range("E41").Value = "  -1.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.99%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.21"
$ws.Range("E43").Value = "  -0.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.822.09"
$ws.Range("E44").Value = "  -1.16%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "55.44"
$ws.Range("E45").Value = "  -3.92%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000104"
$ws.Range("E47").Value = "  -6.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.025"
$ws.Range("E48").Value = "  -1.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05236"
$ws.Range("E49").Value = "  -0.88%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4246"
$ws.Range("E50").Value = "  -1.21%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.976"
$ws.Range("E51").Value = "  -1.53%  "

